$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.03958866666666667
$ws.Range("H2").Value = 0.118766
$ws.Range("I2").Value = 0.0007442768123675562
$ws.Range("J2").Value = 0.0007442768123675561
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 6.655243917690222
$ws.Range("R2").Value = 59.897195259212
$ws.Range("S2").Value = 0.000222105765853226
$ws.Range("T2").Value = 0.000222105765853226

$ws.Range("G3").Value = 0.03958866666666667
$ws.Range("H3").Value = 0.118766
$ws.Range("I3").Value = 0.0007442768123675562
$ws.Range("J3").Value = 0.0007442768123675561
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 6.453199528395778
$ws.Range("R3").Value = 58.078795755562
$ws.Range("S3").Value = 0.0002153629290202577
$ws.Range("T3").Value = 0.0002153629290202577

$ws.Range("G4").Value = 0.03958866666666667
$ws.Range("H4").Value = 0.118766
$ws.Range("I4").Value = 0.0007442768123675562
$ws.Range("J4").Value = 0.0007442768123675561
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 6.571462527993333
$ws.Range("R4").Value = 59.14316275194
$ws.Range("S4").Value = 0.0002193097256249464
$ws.Range("T4").Value = 0.0002193097256249464

$ws.Range("G5").Value = 0.03958866666666667
$ws.Range("H5").Value = 0.118766
$ws.Range("I5").Value = 0.0007442768123675562
$ws.Range("J5").Value = 0.0007442768123675561
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 2.621828109944222
$ws.Range("R5").Value = 23.596452989498
$ws.Range("S5").Value = 0.00008749839186912617
$ws.Range("T5").Value = 0.00008749839186912616

$ws.Range("H6").Value = 0.059669
$ws.Range("I6").Value = 0.0003739306966401135
$ws.Range("J6").Value = 0.0003739306966401134
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 3.343648429050889
$ws.Range("R6").Value = 30.092835861458
$ws.Range("S6").Value = 0.0001115877350647167
$ws.Range("T6").Value = 0.0001115877350647167

$ws.Range("H7").Value = 0.059669
$ws.Range("I7").Value = 0.0003739306966401135
$ws.Range("J7").Value = 0.0003739306966401134
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("Q7").Value = 3.242139691998111
$ws.Range("S7").Value = 0.0001082000792458259
$ws.Range("T7").Value = 0.0001082000792458259

$ws.Range("H8").Value = 0.059669
$ws.Range("I8").Value = 0.0003739306966401135
$ws.Range("J8").Value = 0.0003739306966401134
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 3.301555980523333
$ws.Range("R8").Value = 29.71400382471
$ws.Range("S8").Value = 0.0001101829818156284
$ws.Range("T8").Value = 0.0001101829818156284

$ws.Range("H9").Value = 0.059669
$ws.Range("I9").Value = 0.0003739306966401135
$ws.Range("J9").Value = 0.0003739306966401134
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 1.317227670311889
$ws.Range("R9").Value = 11.855049032807
$ws.Range("S9").Value = 0.00004395990051394246
$ws.Range("T9").Value = 0.00004395990051394245

$ws.Range("G10").Value = 53.131305
$ws.Range("H10").Value = 159.393915
$ws.Range("I10").Value = 0.9988817924909924
$ws.Range("J10").Value = 0.9988817924909923
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 8931.894509544669
$ws.Range("R10").Value = 80387.05058590202
$ws.Range("S10").Value = 0.298084532302334
$ws.Range("T10").Value = 0.2980845323023341

$ws.Range("G11").Value = 53.131305
$ws.Range("H11").Value = 159.393915
$ws.Range("I11").Value = 0.9988817924909924
$ws.Range("J11").Value = 0.9988817924909923
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 8660.734024107545
$ws.Range("R11").Value = 77946.60621696791
$ws.Range("S11").Value = 0.2890350807672734
$ws.Range("T11").Value = 0.2890350807672734

$ws.Range("G12").Value = 53.131305
$ws.Range("H12").Value = 159.393915
$ws.Range("I12").Value = 0.9988817924909924
$ws.Range("J12").Value = 0.9988817924909923
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 8819.452870456649
$ws.Range("R12").Value = 79375.07583410984
$ws.Range("S12").Value = 0.2943320122336024
$ws.Range("T12").Value = 0.2943320122336024

$ws.Range("G13").Value = 53.131305
$ws.Range("H13").Value = 159.393915
$ws.Range("I13").Value = 0.9988817924909924
$ws.Range("J13").Value = 0.9988817924909923
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 3518.712821018305
$ws.Range("R13").Value = 31668.41538916474
$ws.Range("S13").Value = 0.1174301671877826
$ws.Range("T13").Value = 0.1174301671877826
